# Deploy the implementation guide.
# Update the CodeSystem metadata on the "Metadata" sheet:
#   - Status (B6): active -> draft
#   - Date   (B8): 2023-05-12T12:33:13+00:00 -> 2023-08-01T16:12:28+00:00

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2023-08-01T16:12:28+00:00"
